# Tripadvisor New Orleans shard 101 - "update new orleans xlsx files"
#
# Two logical changes:
#   1. hotel_info: insert a new "State" column (value "Louisiana") between
#      "Hotel_Name" and "City" for the header row and the single data row.
#   2. Re-order the sheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Add the "State" column to hotel_info -------------------------------
$hotel = $wb.Worksheets.Item("hotel_info")

# Hotel_Name is column B, City is column C -> insert a fresh column C and
# push City/Zip/... one slot to the right.
$hotel.Columns("C").Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Put review_info ahead of hotel_info in the tab order --------------
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
